$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the numeric-looking "10" value to "21" (kept as text, same as original cell type)
$ws.Range("C2").Value = "21"

# Update the registration date text from 10/11/2021 to 21/12/2021 (kept as text)
$ws.Range("J2").Value = "21/12/2021"

# Update the selected cell in the sheet view to I12 (as reflected in the saved file)
$ws.Range("I12").Select()
